$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin name (B), link (C) and volume% (E) cells are plain text already;
# only the Price column (D) holds numeric-looking text ("1.00", "298.80",
# "23.70", ...) that Excel would otherwise coerce into a Number and so
# lose the exact formatting / trailing zeros. Mark each Price cell we
# touch as Text first so the literal string is preserved.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.045.12'
$ws.Range("E2").Value = '  -1.21%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.268.44'
$ws.Range("E3").Value = '  -1.47%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '298.80'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.23'
$ws.Range("E6").Value = '  -4.28%  '
$ws.Range("E7").Value = '  -2.68%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -2.75%  '
$ws.Range("E10").Value = '  -5.01%  '
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '48.17'
$ws.Range("E12").Value = '  -6.84%  '
$ws.Range("E13").Value = '  +0.98%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.86'
$ws.Range("E14").Value = '  +1.59%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.65'
$ws.Range("E15").Value = '  -1.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.618.96'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.272.35'
$ws.Range("E17").Value = '  -1.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.782'
$ws.Range("E18").Value = '  -2.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.028.51'
$ws.Range("E19").Value = '  -1.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.65'
$ws.Range("E20").Value = '  +1.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0888'
$ws.Range("E21").Value = '  -1.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.96'
$ws.Range("E22").Value = '  -1.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.17'
$ws.Range("E23").Value = '  -2.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '234.96'
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("E25").Value = '  -0.80%  '
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("E27").Value = '  -2.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.70'
$ws.Range("E28").Value = '  -5.01%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '167.45'
$ws.Range("E29").Value = '  +2.27%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.07'
$ws.Range("E30").Value = '  -10.02%  '
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.12'
$ws.Range("E31").Value = '  -0.50%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.53'
$ws.Range("E32").Value = '  -3.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.66'
$ws.Range("E34").Value = '  +5.41%  '
$ws.Range("E35").Value = '  -2.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.36'
$ws.Range("E36").Value = '  -2.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.69'
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0685'
$ws.Range("E38").Value = '  -2.91%  '
$ws.Range("E39").Value = '  -2.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0982'
$ws.Range("E40").Value = '  -1.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.71'
$ws.Range("E42").Value = '  -4.85%  '
$ws.Range("E43").Value = '  -5.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.954.49'
$ws.Range("E44").Value = '  -0.53%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0276'
$ws.Range("E45").Value = '  -1.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.42'
$ws.Range("E46").Value = '  -5.86%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.52'
$ws.Range("E47").Value = '  -6.66%  '
$ws.Range("E48").Value = '  -4.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.491.41'
$ws.Range("E49").Value = '  -1.43%  '
$ws.Range("B50").Value = 'HuobiToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.75'
$ws.Range("E50").Value = '  -3.70%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.75'
$ws.Range("E51").Value = '  -6.91%  '
